$wb = $excel.ActiveWorkbook

# Sheet: u_MAB
$ws1 = $wb.Worksheets.Item("u_MAB")
$ws1.Range("B15").Value = 0
$ws1.Range("A16").Value = 0.2932984991222062
$ws1.Range("B16").Value = 0.05177296069876198
$ws1.Range("B23").Value = 1.069044140854728
$ws1.Range("B24").Value = 0.0372716770134485
$ws1.Range("B25").Value = 0.4930952283364213
$ws1.Range("A27").Value = 0.2086382309069607
$ws1.Range("A40").Value = 0
$ws1.Range("B40").Value = 0
$ws1.Range("A47").Value = 0.3163352131720827
$ws1.Range("A48").Value = 0.2580229554015484
$ws1.Range("B48").Value = 0.313175975823809
$ws1.Range("A49").Value = 0.1730937984422602
$ws1.Range("B61").Value = 0

# Sheet: u_EOH
$ws2 = $wb.Worksheets.Item("u_EOH")
$ws2.Range("A2").Value = -0.3670412353766736
$ws2.Range("A3").Value = -0.2962123418998799

# Sheet: v_l
$ws3 = $wb.Worksheets.Item("v_l")
$ws3.Range("A2").Value = 4027805.565555137
$ws3.Range("A3").Value = 3040539.913340235
$ws3.Range("A4").Value = 0
